$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item(1)

# Add the new "UserEmail" sheet right after "Login"
$userEmailSheet = $wb.Worksheets.Add([System.Type]::Missing, $loginSheet)
$userEmailSheet.Name = "UserEmail"

# Populate UserEmail sheet
$userEmailSheet.Range("A1").Value = "emailid"
$userEmailSheet.Range("A2").Value = "manuav@2003@gmail.com"
$userEmailSheet.Hyperlinks.Add($userEmailSheet.Range("A2"), "mailto:manuav@2003@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "manuav@2003@gmail.com")
$userEmailSheet.Columns.Item(1).AutoFit()
$userEmailSheet.Range("A2").Select()

# Add the email row on the Login sheet
$loginSheet.Range("A4").Value = "manuav@2003@gmail.com"
$loginSheet.Hyperlinks.Add($loginSheet.Range("A4"), "mailto:manuav@2003@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "manuav@2003@gmail.com")

# Restore Login as the active sheet/tab with its own selection
$loginSheet.Activate()
$loginSheet.Range("F10").Select()
